# "Add files via upload" - appends two new daily error-count rows (41 & 42)
# to the "Daily 100 Error Counts" sheet, and leaves the selection on G47.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 41: 12/8/2025 (serial 45999)
$ws.Range("A41").Value = 45999
$ws.Range("B41").Value = 620
$ws.Range("C41").Value = 18
$ws.Range("D41").Value = 602

# Row 42: 12/5/2025 (serial 45996)
$ws.Range("A42").Value = 45996
$ws.Range("B42").Value = 664
$ws.Range("C42").Value = 15
$ws.Range("D42").Value = 649

# Match the saved cursor/selection position recorded in the workbook
$ws.Range("G47").Select()
